$wb = $excel.ActiveWorkbook

# Column F ("想去人数" - number of people interested) updates.
# The same updates apply identically to both the "展览" and "全部类型" sheets.
$updates = @{
    2  = 86
    5  = 2806
    8  = 12
    9  = 1503
    13 = 1251
    15 = 387
    18 = 44
    20 = 82
    22 = 2769
    23 = 333
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
